$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Transcriptions sheet: fill in the new M110 entry (row 16) and the
#    outstanding "Transcription Draft?"/Notes for M109 (row 15).
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Transcriptions")

$tr.Cells.Item(15, 7).Value = "Y"
$tr.Cells.Item(15, 8).Value = "Look at final few folios once the facsimiles arrive or when I get to the archive in person, as the microfilm just doesn't make any sense…"

$tr.Cells.Item(16, 2).Value = "M110 [The Reader]. Typescript fragment, with the author's ms. corrections, unsigned and undated. 5p."
$tr.Cells.Item(16, 3).Value = "m110"
$tr.Cells.Item(16, 4).Value = "transcriptions/m110.xml"
$tr.Cells.Item(16, 6).Value = 5

$tr.Rows.Item(16).RowHeight = 32

$tr.Range("A16").Select()

# ---------------------------------------------------------------------------
# 2. Annotations sheet: re-sort the whole table (A2:F154) alphabetically by
#    the Title column, matching the author's Data > Sort action.
# ---------------------------------------------------------------------------
$an = $wb.Worksheets.Item("Annotations")

$sortRange = $an.Range("A2:F154")
$sortKey = $an.Range("A2")
$sortRange.Sort($sortKey, 1)

$an.Range("A146").Select()
